# CryCompanywiseStockReport_1.xlsx - "Add file from private repo"
#
# The source report lists two rows per stock-code whose Rate/Value/Qty/Amount
# columns (B:G) had been entered against the wrong stock code; this edit
# swaps the B:G contents of each such row-pair so the figures line up with
# the correct code, while the serial numbers in column A stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData($row1, $row2) {
    $rangeA = $ws.Range("B$row1" + ":G$row1")
    $rangeB = $ws.Range("B$row2" + ":G$row2")

    $valA = $rangeA.Value2
    $valB = $rangeB.Value2

    $rangeA.Value = $valB
    $rangeB.Value = $valA
}

# (row, row) pairs whose B:G values are exchanged
Swap-RowData 142 143
Swap-RowData 176 177
Swap-RowData 256 257
Swap-RowData 305 306
Swap-RowData 309 310
Swap-RowData 342 344
Swap-RowData 347 348
Swap-RowData 364 365
Swap-RowData 367 368
Swap-RowData 374 375
Swap-RowData 392 393
Swap-RowData 423 424
Swap-RowData 449 450
Swap-RowData 528 529
Swap-RowData 571 572
Swap-RowData 582 583
Swap-RowData 585 586
Swap-RowData 591 592
Swap-RowData 593 594
Swap-RowData 596 597
Swap-RowData 679 680
Swap-RowData 701 702
Swap-RowData 712 713
